# Script applied 11-11-2023 20:45 — refresh odds/results for matches whose
# home/away rows had been recorded swapped, and append newly played matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Pairs of rows whose match-detail columns (F:V — teams, goals, odds,
#    timestamps, url) were swapped between the two rows. Columns A:E
#    (Indice/pais/torneio/temporada/data_partida) stay put.
# ---------------------------------------------------------------------------
$pairs = @(
    @(4, 5),
    @(8, 9),
    @(26, 27),
    @(36, 37),
    @(48, 49),
    @(53, 54),
    @(86, 87)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = "F" + $r1 + ":V" + $r1
    $range2 = "F" + $r2 + ":V" + $r2
    $tmp = $ws.Range($range1).Value2
    $ws.Range($range1).Value2 = $ws.Range($range2).Value2
    $ws.Range($range2).Value2 = $tmp
}

# ---------------------------------------------------------------------------
# 2) Append three newly finished matches as rows 114-116.
# ---------------------------------------------------------------------------

# Row 114
$ws.Cells.Item(113, 1).Copy()
$ws.Cells.Item(114, 1).PasteSpecial(-4122)
$ws.Cells.Item(113, 5).Copy()
$ws.Cells.Item(114, 5).PasteSpecial(-4122)

$ws.Cells.Item(114, 1).Value2 = 113
$ws.Cells.Item(114, 2).Value2 = "italy"
$ws.Cells.Item(114, 3).Value2 = "serie-a"
$ws.Cells.Item(114, 4).Value2 = "2023-2024"
$ws.Cells.Item(114, 5).Value2 = 45241.625
$ws.Cells.Item(114, 6).Value2 = "Lecce"
$ws.Cells.Item(114, 7).Value2 = 2
$ws.Cells.Item(114, 8).Value2 = "AC Milan"
$ws.Cells.Item(114, 9).Value2 = 2
$ws.Cells.Item(114, 10).Value2 = 4.17
$ws.Cells.Item(114, 11).Value2 = "29/10/2023 11:02"
$ws.Cells.Item(114, 12).Value2 = 5.14
$ws.Cells.Item(114, 13).Value2 = "11/11/2023 14:58"
$ws.Cells.Item(114, 14).Value2 = 3.42
$ws.Cells.Item(114, 15).Value2 = "29/10/2023 11:02"
$ws.Cells.Item(114, 16).Value2 = 3.56
$ws.Cells.Item(114, 17).Value2 = "11/11/2023 14:57"
$ws.Cells.Item(114, 18).Value2 = 1.96
$ws.Cells.Item(114, 19).Value2 = "29/10/2023 11:02"
$ws.Cells.Item(114, 20).Value2 = 1.8
$ws.Cells.Item(114, 21).Value2 = "11/11/2023 14:58"
$ws.Cells.Item(114, 22).Value2 = "https://www.betexplorer.com/football/italy/serie-a/lecce-ac-milan/xYwJ2D2G/"

# Row 115
$ws.Cells.Item(113, 1).Copy()
$ws.Cells.Item(115, 1).PasteSpecial(-4122)
$ws.Cells.Item(113, 5).Copy()
$ws.Cells.Item(115, 5).PasteSpecial(-4122)

$ws.Cells.Item(115, 1).Value2 = 114
$ws.Cells.Item(115, 2).Value2 = "italy"
$ws.Cells.Item(115, 3).Value2 = "serie-a"
$ws.Cells.Item(115, 4).Value2 = "2023-2024"
$ws.Cells.Item(115, 5).Value2 = 45241.75
$ws.Cells.Item(115, 6).Value2 = "Juventus"
$ws.Cells.Item(115, 7).Value2 = 2
$ws.Cells.Item(115, 8).Value2 = "Cagliari"
$ws.Cells.Item(115, 9).Value2 = 1
$ws.Cells.Item(115, 10).Value2 = 1.41
$ws.Cells.Item(115, 11).Value2 = "29/10/2023 11:22"
$ws.Cells.Item(115, 12).Value2 = 1.46
$ws.Cells.Item(115, 13).Value2 = "11/11/2023 17:38"
$ws.Cells.Item(115, 14).Value2 = 4.73
$ws.Cells.Item(115, 15).Value2 = "29/10/2023 11:22"
$ws.Cells.Item(115, 16).Value2 = 4.62
$ws.Cells.Item(115, 17).Value2 = "11/11/2023 17:59"
$ws.Cells.Item(115, 18).Value2 = 8.19
$ws.Cells.Item(115, 19).Value2 = "29/10/2023 11:22"
$ws.Cells.Item(115, 20).Value2 = 8.26
$ws.Cells.Item(115, 21).Value2 = "11/11/2023 17:59"
$ws.Cells.Item(115, 22).Value2 = "https://www.betexplorer.com/football/italy/serie-a/juventus-cagliari/U7zB4ZX3/"

# Row 116
$ws.Cells.Item(113, 1).Copy()
$ws.Cells.Item(116, 1).PasteSpecial(-4122)
$ws.Cells.Item(113, 5).Copy()
$ws.Cells.Item(116, 5).PasteSpecial(-4122)

$ws.Cells.Item(116, 1).Value2 = 115
$ws.Cells.Item(116, 2).Value2 = "italy"
$ws.Cells.Item(116, 3).Value2 = "serie-a"
$ws.Cells.Item(116, 4).Value2 = "2023-2024"
$ws.Cells.Item(116, 5).Value2 = 45241.86458333334
$ws.Cells.Item(116, 6).Value2 = "Monza"
$ws.Cells.Item(116, 7).Value2 = 1
$ws.Cells.Item(116, 8).Value2 = "Torino"
$ws.Cells.Item(116, 9).Value2 = 1
$ws.Cells.Item(116, 10).Value2 = 2.33
$ws.Cells.Item(116, 11).Value2 = "29/10/2023 11:02"
$ws.Cells.Item(116, 12).Value2 = 2.79
$ws.Cells.Item(116, 13).Value2 = "11/11/2023 20:44"
$ws.Cells.Item(116, 14).Value2 = 3.25
$ws.Cells.Item(116, 15).Value2 = "29/10/2023 11:02"
$ws.Cells.Item(116, 16).Value2 = 3.13
$ws.Cells.Item(116, 17).Value2 = "11/11/2023 20:33"
$ws.Cells.Item(116, 18).Value2 = 3.08
$ws.Cells.Item(116, 19).Value2 = "29/10/2023 11:02"
$ws.Cells.Item(116, 20).Value2 = 2.86
$ws.Cells.Item(116, 21).Value2 = "11/11/2023 20:44"
$ws.Cells.Item(116, 22).Value2 = "https://www.betexplorer.com/football/italy/serie-a/monza-torino/GOxN1XHM/"
